$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet 1
$ws1.Range("C2").Value = 0.003769909960995901
$ws1.Range("D2").Value = 0.007793352989302926
$ws1.Range("E2").Value = 0.002113201754037048
$ws1.Range("F2").Value = 0.003702305540019583
$ws1.Range("G2").Value = 0.0003887404917641695
$ws1.Range("I2").Value = 0.003723787788712729
$ws1.Range("J2").Value = 0.007744018133863363
$ws1.Range("K2").Value = 0.002106521197767294
$ws1.Range("L2").Value = 0.003677948519669682
$ws1.Range("M2").Value = 0.0003881208956607774
$ws1.Range("N2").Value = 0.2084769676046654
$ws1.Range("O2").Value = 0.1160767438553549
$ws1.Range("C3").Value = 0.003769909960995901
$ws1.Range("D3").Value = 0.007793352989302926
$ws1.Range("E3").Value = 0.002113201754037048
$ws1.Range("F3").Value = 0.003702305540019583
$ws1.Range("G3").Value = 0.0003887404917641695
$ws1.Range("I3").Value = 0.003720913787828295
$ws1.Range("J3").Value = 0.007734087273310928
$ws1.Range("K3").Value = 0.002102735184248351
$ws1.Range("L3").Value = 0.003673164546190823
$ws1.Range("M3").Value = 0.0003877641643389194
$ws1.Range("N3").Value = 0.173005804283726
$ws1.Range("O3").Value = 0.07992939702573887
$ws1.Range("C4").Value = 0.003771759255797842
$ws1.Range("D4").Value = 0.00778260569127138
$ws1.Range("E4").Value = 0.002109276911031228
$ws1.Range("F4").Value = 0.003697117852169486
$ws1.Range("G4").Value = 0.0003883691171730131
$ws1.Range("I4").Value = 0.003723787788712729
$ws1.Range("J4").Value = 0.007744018133863363
$ws1.Range("K4").Value = 0.002106521197767294
$ws1.Range("L4").Value = 0.003677948519669682
$ws1.Range("M4").Value = 0.0003881208956607774
$ws1.Range("N4").Value = 0.1770877874050666
$ws1.Range("O4").Value = 0.08469559471030816
$ws1.Range("C5").Value = 0.003771759255797842
$ws1.Range("D5").Value = 0.00778260569127138
$ws1.Range("E5").Value = 0.002109276911031228
$ws1.Range("F5").Value = 0.003697117852169486
$ws1.Range("G5").Value = 0.0003883691171730131
$ws1.Range("I5").Value = 0.003720913787828295
$ws1.Range("J5").Value = 0.007734087273310928
$ws1.Range("K5").Value = 0.002102735184248351
$ws1.Range("L5").Value = 0.003673164546190823
$ws1.Range("M5").Value = 0.0003877641643389194
$ws1.Range("N5").Value = 0.141616624084127
$ws1.Range("O5").Value = 0.04854824788069214
$ws1.Range("C6").Value = 0.003775157024594888
$ws1.Range("D6").Value = 0.007787727026405036
$ws1.Range("E6").Value = 0.002108787100421698
$ws1.Range("F6").Value = 0.003699702839221671
$ws1.Range("G6").Value = 0.000388320234647194
$ws1.Range("I6").Value = 0.003725249565827561
$ws1.Range("J6").Value = 0.007742350711306432
$ws1.Range("K6").Value = 0.002105977697171425
$ws1.Range("L6").Value = 0.003677140371788187
$ws1.Range("M6").Value = 0.0003880697129385161
$ws1.Range("N6").Value = 0.1496888401862124
$ws1.Range("O6").Value = 0.0577155186099103
$ws1.Range("C7").Value = 0.003767156216784404
$ws1.Range("D7").Value = 0.007780855852146149
$ws1.Range("E7").Value = 0.002108585252719454
$ws1.Range("F7").Value = 0.003696275723234606
$ws1.Range("G7").Value = 0.0003883035805857967
$ws1.Range("I7").Value = 0.003725249565827561
$ws1.Range("J7").Value = 0.007742350711306432
$ws1.Range("K7").Value = 0.002105977697171425
$ws1.Range("L7").Value = 0.003677140371788187
$ws1.Range("M7").Value = 0.0003880697129385161
$ws1.Range("N7").Value = 0.2625447922510269
$ws1.Range("O7").Value = 0.1705425960019301
$ws1.Range("C8").Value = 0.003774011943783199
$ws1.Range("D8").Value = 0.007787067112526803
$ws1.Range("E8").Value = 0.002110763953587128
$ws1.Range("F8").Value = 0.003699278179158993
$ws1.Range("G8").Value = 0.0003885097115356922
$ws1.Range("I8").Value = 0.003725249565827561
$ws1.Range("J8").Value = 0.007742350711306432
$ws1.Range("K8").Value = 0.002106521197767294
$ws1.Range("L8").Value = 0.003677948519669682
$ws1.Range("M8").Value = 0.0003881208956607774
$ws1.Range("N8").Value = 0.1685526780378239
$ws1.Range("O8").Value = 0.0765745300688003

# Sheet 2
$ws2.Range("C2").Value = 0.003769909960995901
$ws2.Range("D2").Value = 0.007793352989302926
$ws2.Range("E2").Value = 0.002113201754037048
$ws2.Range("F2").Value = 0.003702305540019583
$ws2.Range("G2").Value = 0.0003887404917641695
$ws2.Range("I2").Value = 0.003723787788712729
$ws2.Range("J2").Value = 0.007744018133863363
$ws2.Range("K2").Value = 0.002106521197767294
$ws2.Range("L2").Value = 0.003677948519669682
$ws2.Range("M2").Value = 0.0003881208956607774
$ws2.Range("N2").Value = 0.2084769676046654
$ws2.Range("O2").Value = 0.1160767438553549
$ws2.Range("C3").Value = 0.003769909960995901
$ws2.Range("D3").Value = 0.007793352989302926
$ws2.Range("E3").Value = 0.002113201754037048
$ws2.Range("F3").Value = 0.003702305540019583
$ws2.Range("G3").Value = 0.0003887404917641695
$ws2.Range("I3").Value = 0.003720913787828295
$ws2.Range("J3").Value = 0.007734087273310928
$ws2.Range("K3").Value = 0.002102735184248351
$ws2.Range("L3").Value = 0.003673164546190823
$ws2.Range("M3").Value = 0.0003877641643389194
$ws2.Range("N3").Value = 0.173005804283726
$ws2.Range("O3").Value = 0.07992939702573887
$ws2.Range("C4").Value = 0.003771759255797842
$ws2.Range("D4").Value = 0.00778260569127138
$ws2.Range("E4").Value = 0.002109276911031228
$ws2.Range("F4").Value = 0.003697117852169486
$ws2.Range("G4").Value = 0.0003883691171730131
$ws2.Range("I4").Value = 0.003723787788712729
$ws2.Range("J4").Value = 0.007744018133863363
$ws2.Range("K4").Value = 0.002106521197767294
$ws2.Range("L4").Value = 0.003677948519669682
$ws2.Range("M4").Value = 0.0003881208956607774
$ws2.Range("N4").Value = 0.1770877874050666
$ws2.Range("O4").Value = 0.08469559471030816
$ws2.Range("C5").Value = 0.003771759255797842
$ws2.Range("D5").Value = 0.00778260569127138
$ws2.Range("E5").Value = 0.002109276911031228
$ws2.Range("F5").Value = 0.003697117852169486
$ws2.Range("G5").Value = 0.0003883691171730131
$ws2.Range("I5").Value = 0.003720913787828295
$ws2.Range("J5").Value = 0.007734087273310928
$ws2.Range("K5").Value = 0.002102735184248351
$ws2.Range("L5").Value = 0.003673164546190823
$ws2.Range("M5").Value = 0.0003877641643389194
$ws2.Range("N5").Value = 0.141616624084127
$ws2.Range("O5").Value = 0.04854824788069214
$ws2.Range("C6").Value = 0.003775157024594888
$ws2.Range("D6").Value = 0.007787727026405036
$ws2.Range("E6").Value = 0.002108787100421698
$ws2.Range("F6").Value = 0.003699702839221671
$ws2.Range("G6").Value = 0.000388320234647194
$ws2.Range("I6").Value = 0.003725249565827561
$ws2.Range("J6").Value = 0.007742350711306432
$ws2.Range("K6").Value = 0.002105977697171425
$ws2.Range("L6").Value = 0.003677140371788187
$ws2.Range("M6").Value = 0.0003880697129385161
$ws2.Range("N6").Value = 0.1496888401862124
$ws2.Range("O6").Value = 0.0577155186099103
$ws2.Range("C7").Value = 0.003767156216784404
$ws2.Range("D7").Value = 0.007780855852146149
$ws2.Range("E7").Value = 0.002108585252719454
$ws2.Range("F7").Value = 0.003696275723234606
$ws2.Range("G7").Value = 0.0003883035805857967
$ws2.Range("I7").Value = 0.003725249565827561
$ws2.Range("J7").Value = 0.007742350711306432
$ws2.Range("K7").Value = 0.002105977697171425
$ws2.Range("L7").Value = 0.003677140371788187
$ws2.Range("M7").Value = 0.0003880697129385161
$ws2.Range("N7").Value = 0.2625447922510269
$ws2.Range("O7").Value = 0.1705425960019301
$ws2.Range("C8").Value = 0.003774011943783199
$ws2.Range("D8").Value = 0.007787067112526803
$ws2.Range("E8").Value = 0.002110763953587128
$ws2.Range("F8").Value = 0.003699278179158993
$ws2.Range("G8").Value = 0.0003885097115356922
$ws2.Range("I8").Value = 0.003725249565827561
$ws2.Range("J8").Value = 0.007742350711306432
$ws2.Range("K8").Value = 0.002106521197767294
$ws2.Range("L8").Value = 0.003677948519669682
$ws2.Range("M8").Value = 0.0003880697129385161
$ws2.Range("N8").Value = 0.1685526780378239
$ws2.Range("O8").Value = 0.0765745300688003
